$d = $word.ActiveDocument

# 1) Intro / role / company paragraph
$d.Content.Find.Execute(
    "Nous soussignés, QuantFactory, attestons que Monsieur **Mouad Med** occupe au sein de notre entreprise le poste de **Full-stack developer** au sein du département **IT**, et ce depuis son embauche en **Contrat à Durée Indéterminée (CDI)** le **21 mai 2023**.  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Je soussigné(e), [Votre Nom], représentant(e) des Ressources Humaines chez QuantFactory, atteste par la présente que Monsieur **Mouad Mo** occupe le poste de **IA Developer** au sein de notre entreprise.  ",
    2) | Out-Null

# 2) Employment / hire-date / department paragraph
$d.Content.Find.Execute(
    "À ce jour, Monsieur **Mouad Med** exerce toujours ses fonctions avec professionnalisme et dévouement.  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Monsieur Mo a été embauché(e) en **Contrat à Durée Indéterminée (CDI)** le **30 mai 2023** et exerce ses fonctions au sein du département **IT**.  ",
    2) | Out-Null

# 3) Purpose of the attestation paragraph
$d.Content.Find.Execute(
    "La présente attestation est délivrée à sa demande pour servir et valoir ce que de droit.  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cette attestation est délivrée à sa demande pour les usages qu’il jugera nécessaires.  ",
    2) | Out-Null

# 4) Date line -> signatory name
$d.Content.Find.Execute(
    "**Le [27/05/2025]**  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Votre Nom]  ",
    2) | Out-Null

# 5) "Pour QuantFactory," -> role title
$d.Content.Find.Execute(
    "**Pour QuantFactory,**  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Responsable des Ressources Humaines  ",
    2) | Out-Null

# 6) [Signature] -> QuantFactory
$d.Content.Find.Execute(
    "[Signature]  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "QuantFactory  ",
    2) | Out-Null

# 7) [Nom du responsable RH] -> [Coordonnées de l'entreprise]
$d.Content.Find.Execute(
    "[Nom du responsable RH]  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Coordonnées de l’entreprise]  ",
    2) | Out-Null

# 8) Service des Ressources Humaines -> Le 30 mai 2025
$d.Content.Find.Execute(
    "**Service des Ressources Humaines**",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Le **30 mai 2025**",
    2) | Out-Null
